$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "13 Nov 2025, 09:20 AM"

# --- Stock List sheet: refresh the stock rows (2 rows dropped from the
#     top, two new rows appended at the bottom, everything else shifts) ---
$ws = $wb.Worksheets.Item("Stock List")

$data = @(
    ("SMLMAH", 3018.2, 0.0464, 0),
    ("LENSKART", 421.2, 0, 73072.7323),
    ("STUDDS", 549.55, -0.0182, 2163.0596),
    ("PIRAMALFIN", 1415, -1.2699, 0),
    ("ORKLAINDIA", 667, -0.1123, 9147.4558),
    ("MIDWESTLTD", 1130, 2.9331, 3969.7344),
    ("CANHLIFE", 121, -0.3213, 11532.05),
    ("RUBICON", 625, -0.3587, 10333.9631),
    ("CRAMC", 312.2, -0.2556, 6241.7655),
    ("LGEINDIA", 1660, 0.1629, 112492.9485),
    ("TATACAP", 325, -0.0615, 138043.1411),
    ("WEWORK", 608, -0.7428, 8209.5947),
    ("ADVANCE", 142.99, 0.598, 913.7572),
    ("OMFREIGHT", 90.99, 1.0663, 303.18),
    ("GLOTTIS", 72.5, -0.0551, 670.2914),
    ("FABTECH", 229.4, 0.1309, 1018.3677),
    ("PACEDIGITK", 217.9, 0.0873, 4699.2992),
    ("JAINREC", 426, 0.3415, 14650.6182),
    ("EPACKPEB", 314.6, 1.3041, 3119.5368),
    ("BMWVENTLTD", 65.16, -0.7615, 569.3707),
    ("STYL", 351.4, 1.0932, 5624.3705),
    ("JARO", 626.15, -0.4847, 1394.0737),
    ("SOLARWORLD", 291, 0.0172, 2521.7409),
    ("ARSSBL", 667.3, 2.4252, 4086.55),
    ("GANESHCP", 287.8, 1.1955, 1149.3442),
    ("ATLANTAELE", 1025, 1.3547, 7776.5557),
    ("GKENERGY", 196.2, -0.1018, 3983.3311),
    ("SAATVIKGL", 467.6, 0, 5943.43),
    ("IVALUE", 320, 0.0469, 1712.4731),
    ("VMSTMT", 68, 0.3394, 336.3507),
    ("EUROPRATIK", 344.8, 0.6275, 3501.883),
    ("SHRINGARMS", 219, 0.1189, 2109.3553),
    ("DEVX", 43.47, -0.1149, 392.4961),
    ("URBANCO", 142.71, -0.007, 20493.1924),
    ("AMANTA", 121.76, -0.8146, 476.6691),
    ("CPEDU", 290.2, -3.3955, 546.5159),
    ("AHCL", 141.4, 0, 751.5622),
    ("STLNETWORK", 23.8, 0.4643, 1155.8851),
    ("VIKRAN", 105.8, 0.3224, 2719.9286),
    ("MEIL", 460, 0.4038, 1265.8741),
    ("SHREEJISPG", 295.1, 1.671, 4728.7023),
    ("GEMAROMA", 205.95, 0.4291, 1071.227),
    ("PATELRMART", 222.09, 0, 741.7923),
    ("VIKRAMSOLR", 316.3, 0, 11441.1097),
    ("REGAAL", 93.22, 4.4833, 916.4987),
    ("BLUESTONE", 596.75, 0.8705, 8952.1132),
    ("ALLTIME", 280.9, 0, 1840.1135),
    ("JSWCEMENT", 122.24, 0.0082, 16664.4096),
    ("HILINFRA", 70.06, 0, 0),
    ("LOTUSDEV", 172.75, 1.118, 8349.3559),
    ("MBEL", 435.5, -0.4002, 2498.8057),
    ("LAXMIINDIA", 151.9, 0.1979, 792.381),
    ("CPPLUS", 1570, 4.9746, 17531.705),
    ("SHANTIGOLD", 228.7, 0.0175, 1648.5471),
    ("BRIGHOTEL", 80.1, 0.125, 3038.7405),
    ("INDIQUBE", 221.98, 0.2031, 4652.3849),
    ("EBGNG", 339.8, 0.4879, 3855.2994),
    ("CHEMBONDCH", 149.9, 0.5568, 400.9473),
    ("ANTHEM", 690, 0.5245, 38548.9139),
    ("SMARTWORKS", 584.8, -0.0171, 6682.7414),
    ("TRAVELFOOD", 1254, -1.6702, 16793.0846),
    ("CRIZAC", 295.95, 0.3731, 5159.359),
    ("IGCL", 100.45, -1.2097, 642.8616),
    ("SAMBHV", 111.5, 0.2067, 3278.809),
    ("HDBFS", 727.2, 0.186, 60214.1009),
    ("ELLEN", 442, -0.0565, 0),
    ("KALPATARU", 374.4, 1.4084, 7602.3446),
    ("GLOBECIVIL", 68.5, 1.617, 402.5658),
    ("RAYMONDREL", 570.15, 0, 3795.7013),
    ("ARISINFRA", 155.55, -0.0514, 1261.3582),
    ("ABLBL", 129.75, -0.1616, 15858.9509),
    ("OSWALPUMPS", 643.95, 1.1307, 7257.5118),
    ("ENRIN", 3283.6001, 0, 116935.729),
    ("BLUSPRING", 77.15, -1.4813, 1166.4229),
    ("DIGITIDE", 145.25, 0.0965, 2161.4049)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ticker = $row[0]
    $price = $row[1]
    $pctChange = $row[2]
    $marketCap = $row[3]
    $ws.Cells.Item($r, 2).Value = $ticker
    $ws.Cells.Item($r, 3).Value = $ticker
    $ws.Cells.Item($r, 4).Value = $price
    $ws.Cells.Item($r, 5).Value = $pctChange
    $ws.Cells.Item($r, 8).Value = $marketCap
}
